# SCD0278 - add two new scripted columns (TEXT5 / TEXT6 headers) and move the
# "Cuti Sakit" answer that used to live under column O (TEXT4) out to the new
# TEXT6 column (Q), leaving column O's row-2 cell empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for SCD0279 - SCD0281 scripting work.
$ws.Range("P1").Value = "TEXT5"
$ws.Range("Q1").Value = "TEXT6"

# Carry the existing O2 answer (with its formatting) over to Q2, then clear O2
# completely (value + formatting) so it no longer occupies that cell.
$ws.Range("O2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)
$ws.Range("Q2").Value = $ws.Range("O2").Text
$ws.Range("O2").Clear()

# Column O no longer needs to be wide now that its wrapped content moved to Q;
# narrow it down to fit the short header text.
$ws.Columns.Item(15).ColumnWidth = 5.5

# Row 2 shrinks now that column O's wrapped paragraph isn't there any more.
$ws.Rows.Item(2).RowHeight = 89.25

# The active selection moves to the relocated answer cell.
[void]$ws.Range("Q2").Select()
